# Auto-generated COM-interop script applying the Bahamut_Profits market-data refresh.
$wb = $excel.ActiveWorkbook

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 933.6
$ws.Range("I125").Value = 810.6667
$ws.Range("J125").Value = 1118
$ws.Range("K125").Value = 7296.0003
$ws.Range("L125").Value = 10062
$ws.Range("M125").Value = -4836.0003
$ws.Range("N125").Value = -14982

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 808.8333
$ws.Range("I137").Value = 760.6923
$ws.Range("J137").Value = 836.04346
$ws.Range("K137").Value = 2282.0769
$ws.Range("L137").Value = 2508.13038
$ws.Range("M137").Value = 267.9231
$ws.Range("N137").Value = -7608.130380000001

# ARM!row32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5385.72
$ws.Range("I32").Value = 4367.42
$ws.Range("J32").Value = 9726.895
$ws.Range("K32").Value = 4367.42
$ws.Range("L32").Value = 9726.895
$ws.Range("M32").Value = -4080.42
$ws.Range("N32").Value = -10300.895

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 970.1905
$ws.Range("I74").Value = 959.1579
$ws.Range("J74").Value = 1075
$ws.Range("K74").Value = 959.1579
$ws.Range("L74").Value = 1075
$ws.Range("M74").Value = -85.15790000000004
$ws.Range("N74").Value = -2823

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 970.1905
$ws.Range("I77").Value = 959.1579
$ws.Range("J77").Value = 1075
$ws.Range("K77").Value = 4795.7895
$ws.Range("L77").Value = 5375
$ws.Range("M77").Value = -427.7894999999999
$ws.Range("N77").Value = -14111

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1255.6061
$ws.Range("I132").Value = 1148.4073
$ws.Range("J132").Value = 1738
$ws.Range("K132").Value = 3445.2219
$ws.Range("L132").Value = 5214
$ws.Range("M132").Value = -915.2219000000005
$ws.Range("N132").Value = -10274

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2474.2263
$ws.Range("I31").Value = 2288.7104
$ws.Range("J31").Value = 2944.2
$ws.Range("K31").Value = 2288.7104
$ws.Range("L31").Value = 2944.2
$ws.Range("M31").Value = -1993.7104
$ws.Range("N31").Value = -3534.2

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2474.2263
$ws.Range("I34").Value = 2288.7104
$ws.Range("J34").Value = 2944.2
$ws.Range("K34").Value = 2288.7104
$ws.Range("L34").Value = 2944.2
$ws.Range("M34").Value = -2086.7104
$ws.Range("N34").Value = -3348.2

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2066
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2066
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6198
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -11098

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1057.2676
$ws.Range("I134").Value = 964.2
$ws.Range("J134").Value = 1564.909
$ws.Range("K134").Value = 2892.6
$ws.Range("L134").Value = 4694.727000000001
$ws.Range("M134").Value = -357.6000000000004
$ws.Range("N134").Value = -9764.727000000001

# CUL!row92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 400
$ws.Range("I92").Value = 300
$ws.Range("J92").Value = 500
$ws.Range("K92").Value = 900
$ws.Range("L92").Value = 1500
$ws.Range("M92").Value = 348
$ws.Range("N92").Value = -3996

# CUL!row103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1263.3334
$ws.Range("I103").Value = 234
$ws.Range("J103").Value = 2550
$ws.Range("K103").Value = 702
$ws.Range("L103").Value = 7650
$ws.Range("M103").Value = 177
$ws.Range("N103").Value = -9408

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 21530.408
$ws.Range("J131").Value = 989.3182
$ws.Range("L131").Value = 2967.9546
$ws.Range("N131").Value = -13047.9546

# CUL!row132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 682.1539
$ws.Range("I132").Value = 607.25
$ws.Range("J132").Value = 802
$ws.Range("K132").Value = 5465.25
$ws.Range("L132").Value = 7218
$ws.Range("M132").Value = -2935.25
$ws.Range("N132").Value = -12278

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4016.5386
$ws.Range("I70").Value = 3880.5264
$ws.Range("K70").Value = 3880.5264
$ws.Range("M70").Value = -3610.5264

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4016.5386
$ws.Range("I73").Value = 3880.5264
$ws.Range("K73").Value = 3880.5264
$ws.Range("M73").Value = -2944.5264

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1616.4445
$ws.Range("I102").Value = 1643.5
$ws.Range("K102").Value = 1643.5
$ws.Range("M102").Value = -21.5

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12983160
$ws.Range("I122").Value = 13304701
$ws.Range("J122").Value = 12500850
$ws.Range("K122").Value = 39914103
$ws.Range("L122").Value = 37502550
$ws.Range("M122").Value = -39911653
$ws.Range("N122").Value = -37507450

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2210.6
$ws.Range("I126").Value = 2083.3809
$ws.Range("J126").Value = 2878.5
$ws.Range("K126").Value = 6250.1427
$ws.Range("L126").Value = 8635.5
$ws.Range("M126").Value = -3780.1427
$ws.Range("N126").Value = -13575.5

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2204.8647
$ws.Range("I132").Value = 2068.6296
$ws.Range("J132").Value = 2572.7
$ws.Range("K132").Value = 6205.888800000001
$ws.Range("L132").Value = 7718.099999999999
$ws.Range("M132").Value = -3675.888800000001
$ws.Range("N132").Value = -12778.1

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4833864.5
$ws.Range("I7").Value = 3062.4119
$ws.Range("J7").Value = 18521136
$ws.Range("K7").Value = 3062.4119
$ws.Range("L7").Value = 18521136
$ws.Range("M7").Value = -2950.4119
$ws.Range("N7").Value = -18521360

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 562501.5
$ws.Range("I40").Value = 674568.1
$ws.Range("K40").Value = 674568.1
$ws.Range("M40").Value = -674432.1

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 4833864.5
$ws.Range("I126").Value = 3062.4119
$ws.Range("J126").Value = 18521136
$ws.Range("K126").Value = 9187.235700000001
$ws.Range("L126").Value = 55563408
$ws.Range("M126").Value = -6717.235700000001
$ws.Range("N126").Value = -55568348

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1745.6459
$ws.Range("I132").Value = 1665.5682
$ws.Range("J132").Value = 2626.5
$ws.Range("K132").Value = 4996.7046
$ws.Range("L132").Value = 7879.5
$ws.Range("M132").Value = -2466.7046
$ws.Range("N132").Value = -12939.5

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 917.5454999999999
$ws.Range("I126").Value = 773.25
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 2319.75
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 150.25
$ws.Range("N126").Value = -7940

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 526.4286
$ws.Range("I132").Value = 505.81033
$ws.Range("J132").Value = 626.0833
$ws.Range("K132").Value = 1517.43099
$ws.Range("L132").Value = 1878.2499
$ws.Range("M132").Value = 1012.56901
$ws.Range("N132").Value = -6938.2499

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 820.25
$ws.Range("I136").Value = 949.3929000000001
$ws.Range("J136").Value = 594.25
$ws.Range("K136").Value = 2848.1787
$ws.Range("L136").Value = 1782.75
$ws.Range("M136").Value = -298.1787000000004
$ws.Range("N136").Value = -6882.75
